# Generate Report for Handoff
# Adds two new file rows (8d94f427-... and e478cf74-...) to the
# Overview / zh-cn / de-de sheets of the localization-status report,
# mirroring the existing "Ready for handoff" rows already present.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New file identifiers (new rows appended as r=6 and r=7 on every sheet)
# ---------------------------------------------------------------------
$guid1 = "8d94f427-c817-4b74-878b-b6b6dfbbc76d"
$guid2 = "e478cf74-277d-4ce7-bfd4-8f3f2ceade18"

$hash1 = "58f92d8b750e56f12a24b0be191d1ac2714c7234"
$hash2 = "235bd09204594b439ee1a6a6bb3b3b66249b02f8"

$status = "Ready for handoff"

$dateOverview = "2016-17-18 04:17:28"
$dateZh       = "2016-03-18 04:17:24"
$dateDe       = "2016-03-18 04:17:28"

$commitE2e1 = "c0ffee1111111111111111111111111111111a1"
$commitE2e2 = "c0ffee2222222222222222222222222222222a2"
$commitZh1  = "beefcafe111111111111111111111111111111"
$commitZh2  = "beefcafe222222222222222222222222222222"
$commitDe1  = "dadfeed111111111111111111111111111111a"
$commitDe2  = "dadfeed222222222222222222222222222222b"

# =======================================================================
# Sheet 1: Overview
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(6, 2).Value = $status
$wsOverview.Cells.Item(6, 3).Value = $status
$wsOverview.Cells.Item(6, 4).Value = $dateOverview

$wsOverview.Cells.Item(7, 2).Value = $status
$wsOverview.Cells.Item(7, 3).Value = $status
$wsOverview.Cells.Item(7, 4).Value = $dateOverview

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e1/e2e/$guid1.md",
    "",
    "",
    "$guid1.md"
) | Out-Null

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e2/e2e/$guid2.md",
    "",
    "",
    "$guid2.md"
) | Out-Null

# =======================================================================
# Sheet 2: zh-cn
# =======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(6, 3).Value = $status
$wsZh.Cells.Item(6, 5).Value = $dateZh
$wsZh.Cells.Item(6, 8).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(6, 9).Value = "Include"

$wsZh.Cells.Item(7, 3).Value = $status
$wsZh.Cells.Item(7, 5).Value = $dateZh
$wsZh.Cells.Item(7, 8).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(7, 9).Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e1/e2e/$guid1.md",
    "",
    "",
    "$guid1.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e1/e2e/$guid1.md",
    "",
    "",
    ".md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D6"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid1.$hash1.zh-cn.xlf",
    "",
    "",
    "$guid1.$hash1.zh-cn.xlf"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e2/e2e/$guid2.md",
    "",
    "",
    "$guid2.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e2/e2e/$guid2.md",
    "",
    "",
    ".md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D7"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid2.$hash2.zh-cn.xlf",
    "",
    "",
    "$guid2.$hash2.zh-cn.xlf"
) | Out-Null

# =======================================================================
# Sheet 3: de-de
# =======================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(6, 3).Value = $status
$wsDe.Cells.Item(6, 5).Value = $dateDe
$wsDe.Cells.Item(6, 8).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(6, 9).Value = "Include"

$wsDe.Cells.Item(7, 3).Value = $status
$wsDe.Cells.Item(7, 5).Value = $dateDe
$wsDe.Cells.Item(7, 8).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(7, 9).Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e1/e2e/$guid1.md",
    "",
    "",
    "$guid1.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e1/e2e/$guid1.md",
    "",
    "",
    ".md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D6"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid1.$hash1.de-de.xlf",
    "",
    "",
    "$guid1.$hash1.de-de.xlf"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e2/e2e/$guid2.md",
    "",
    "",
    "$guid2.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e2/e2e/$guid2.md",
    "",
    "",
    ".md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D7"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid2.$hash2.de-de.xlf",
    "",
    "",
    "$guid2.$hash2.de-de.xlf"
) | Out-Null

Write-Host "Handoff rows appended for $guid1 and $guid2"
